$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data row (for variety "Super Queen") was inserted into the daily log
# right before the existing row 315, pushing all following rows (old 315-407)
# down by one (new rows 316-408).
$ws.Rows.Item(315).Insert()

# Populate the newly inserted row 315 with its data.
$ws.Cells.Item(315, 1).Value = 5
$ws.Cells.Item(315, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(315, 3).Value = "Maule"
$ws.Cells.Item(315, 4).Value = 44559
$ws.Cells.Item(315, 5).Value = 7
$ws.Cells.Item(315, 6).Value = "Fruta"
$ws.Cells.Item(315, 7).Value = 100103
$ws.Cells.Item(315, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(315, 9).Value = 100103006
$ws.Cells.Item(315, 10).Value = "Nectarín"
$ws.Cells.Item(315, 11).Value = "Super Queen"
$ws.Cells.Item(315, 12).Value = "Primera"
$ws.Cells.Item(315, 13).Value = 500
$ws.Cells.Item(315, 14).Value = 10000
$ws.Cells.Item(315, 15).Value = 10000
$ws.Cells.Item(315, 16).Value = 10000
$ws.Cells.Item(315, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(315, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(315, 19).Value = 556
$ws.Cells.Item(315, 20).Value = 18
